$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.102.93'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '1.787.99'
$ws.Range('E3').Value = '  -2.74%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''224.55'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''32.97'
$ws.Range('E8').Value = '  +2.80%  '
$ws.Range('E9').Value = '  -2.44%  '
$ws.Range('D10').Value = '''0.0710'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').Value = '''0.0931'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '2.045.25'
$ws.Range('E12').Value = '  -2.80%  '
$ws.Range('D13').Value = '1.789.90'
$ws.Range('E13').Value = '  -2.60%  '
$ws.Range('D14').Value = '''10.83'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('D16').Value = '34.067.82'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('E17').Value = '  -4.60%  '
$ws.Range('D18').Value = '''67.87'
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').Value = '''245.33'
$ws.Range('E19').Value = '  -2.49%  '
$ws.Range('D20').Value = '0.0₃0791'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '''0.998'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '''10.82'
$ws.Range('E22').Value = '  -4.51%  '
$ws.Range('E23').Value = '  -4.20%  '
$ws.Range('E24').Value = '  -2.88%  '
$ws.Range('D25').Value = '''160.56'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('D26').Value = '''16.34'
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('E27').Value = '  -2.63%  '
$ws.Range('D28').Value = '''0.112'
$ws.Range('E28').Value = '  -2.96%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('E32').Value = '  -4.09%  '
$ws.Range('D33').Value = '''3.51'
$ws.Range('E33').Value = '  -2.87%  '
$ws.Range('E34').Value = '  -5.63%  '
$ws.Range('D35').Value = '1.398.52'
$ws.Range('E35').Value = '  -4.16%  '
$ws.Range('D36').Value = '''0.642'
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('E38').Value = '  -3.46%  '
$ws.Range('E39').Value = '  +3.10%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = '''0.916'
$ws.Range('E42').Value = '  -2.95%  '
$ws.Range('D43').Value = '''78.18'
$ws.Range('E43').Value = '  -5.15%  '
$ws.Range('D44').Value = '0.0₆0146'
$ws.Range('E44').Value = '  +16.44%  '
$ws.Range('D45').Value = '''1.08'
$ws.Range('E45').Value = '  +1.69%  '
$ws.Range('D46').Value = '''0.0499'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '''12.48'
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''108.02'
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('D49').Value = '''5.88'
$ws.Range('E49').Value = '  -3.63%  '
$ws.Range('D50').Value = '1.945.42'
$ws.Range('E50').Value = '  -2.78%  '
$ws.Range('E51').Value = '  -0.40%  '
